$d = $word.ActiveDocument

$d.Content.Find.Execute("8+69=", $true, $false, $false, $false, $false, $true, 1, $false, "85-77=", 2) | Out-Null
$d.Content.Find.Execute("46-9=", $true, $false, $false, $false, $false, $true, 1, $false, "45+29=", 2) | Out-Null
$d.Content.Find.Execute("79+16=", $true, $false, $false, $false, $false, $true, 1, $false, "7+26=", 2) | Out-Null
$d.Content.Find.Execute("28+15=", $true, $false, $false, $false, $false, $true, 1, $false, "27-18=", 2) | Out-Null
$d.Content.Find.Execute("35+47=", $true, $false, $false, $false, $false, $true, 1, $false, "39+5=", 2) | Out-Null
$d.Content.Find.Execute("37+57=", $true, $false, $false, $false, $false, $true, 1, $false, "6+27=", 2) | Out-Null
$d.Content.Find.Execute("49+18=", $true, $false, $false, $false, $false, $true, 1, $false, "59+17=", 2) | Out-Null
$d.Content.Find.Execute("62-28=", $true, $false, $false, $false, $false, $true, 1, $false, "41-35=", 2) | Out-Null
$d.Content.Find.Execute("15+67=", $true, $false, $false, $false, $false, $true, 1, $false, "27+17=", 2) | Out-Null
$d.Content.Find.Execute("92-76=", $true, $false, $false, $false, $false, $true, 1, $false, "91-62=", 2) | Out-Null
$d.Content.Find.Execute("38+57=", $true, $false, $false, $false, $false, $true, 1, $false, "68+6=", 2) | Out-Null
$d.Content.Find.Execute("67-8=", $true, $false, $false, $false, $false, $true, 1, $false, "84-68=", 2) | Out-Null
$d.Content.Find.Execute("82-25=", $true, $false, $false, $false, $false, $true, 1, $false, "85-9=", 2) | Out-Null
$d.Content.Find.Execute("77+6=", $true, $false, $false, $false, $false, $true, 1, $false, "58-39=", 2) | Out-Null
$d.Content.Find.Execute("47+19=", $true, $false, $false, $false, $false, $true, 1, $false, "75-19=", 2) | Out-Null
$d.Content.Find.Execute("54-27=", $true, $false, $false, $false, $false, $true, 1, $false, "62-17=", 2) | Out-Null
$d.Content.Find.Execute("56+39=", $true, $false, $false, $false, $false, $true, 1, $false, "84-35=", 2) | Out-Null
$d.Content.Find.Execute("47+44=", $true, $false, $false, $false, $false, $true, 1, $false, "81-78=", 2) | Out-Null
$d.Content.Find.Execute("78+5=", $true, $false, $false, $false, $false, $true, 1, $false, "51-7=", 2) | Out-Null
$d.Content.Find.Execute("32-23=", $true, $false, $false, $false, $false, $true, 1, $false, "14+78=", 2) | Out-Null
$d.Content.Find.Execute("74-49=", $true, $false, $false, $false, $false, $true, 1, $false, "51-14=", 2) | Out-Null
$d.Content.Find.Execute("58+35=", $true, $false, $false, $false, $false, $true, 1, $false, "36-28=", 2) | Out-Null
$d.Content.Find.Execute("17+17=", $true, $false, $false, $false, $false, $true, 1, $false, "13+9=", 2) | Out-Null
$d.Content.Find.Execute("32+49=", $true, $false, $false, $false, $false, $true, 1, $false, "88-19=", 2) | Out-Null
$d.Content.Find.Execute("17+45=", $true, $false, $false, $false, $false, $true, 1, $false, "81-48=", 2) | Out-Null
$d.Content.Find.Execute("17+28=", $true, $false, $false, $false, $false, $true, 1, $false, "24+38=", 2) | Out-Null
$d.Content.Find.Execute("65+8=", $true, $false, $false, $false, $false, $true, 1, $false, "68+29=", 2) | Out-Null
$d.Content.Find.Execute("75-39=", $true, $false, $false, $false, $false, $true, 1, $false, "76+18=", 2) | Out-Null
$d.Content.Find.Execute("16+76=", $true, $false, $false, $false, $false, $true, 1, $false, "28+65=", 2) | Out-Null
$d.Content.Find.Execute("49+45=", $true, $false, $false, $false, $false, $true, 1, $false, "55+27=", 2) | Out-Null
$d.Content.Find.Execute("34+9=", $true, $false, $false, $false, $false, $true, 1, $false, "66-47=", 2) | Out-Null
$d.Content.Find.Execute("51-27=", $true, $false, $false, $false, $false, $true, 1, $false, "83-39=", 2) | Out-Null
$d.Content.Find.Execute("72-4=", $true, $false, $false, $false, $false, $true, 1, $false, "33-8=", 2) | Out-Null
$d.Content.Find.Execute("73-46=", $true, $false, $false, $false, $false, $true, 1, $false, "28+7=", 2) | Out-Null
$d.Content.Find.Execute("90-15=", $true, $false, $false, $false, $false, $true, 1, $false, "29+5=", 2) | Out-Null
$d.Content.Find.Execute("6+19=", $true, $false, $false, $false, $false, $true, 1, $false, "19+67=", 2) | Out-Null
$d.Content.Find.Execute("45+49=", $true, $false, $false, $false, $false, $true, 1, $false, "43-36=", 2) | Out-Null
$d.Content.Find.Execute("34+59=", $true, $false, $false, $false, $false, $true, 1, $false, "60-42=", 2) | Out-Null
$d.Content.Find.Execute("72-55=", $true, $false, $false, $false, $false, $true, 1, $false, "58+18=", 2) | Out-Null
$d.Content.Find.Execute("28+59=", $true, $false, $false, $false, $false, $true, 1, $false, "63+8=", 2) | Out-Null
$d.Content.Find.Execute("55+37=", $true, $false, $false, $false, $false, $true, 1, $false, "45-9=", 2) | Out-Null
$d.Content.Find.Execute("72-54=", $true, $false, $false, $false, $false, $true, 1, $false, "48+46=", 2) | Out-Null
$d.Content.Find.Execute("81-27=", $true, $false, $false, $false, $false, $true, 1, $false, "7+59=", 2) | Out-Null
$d.Content.Find.Execute("81-42=", $true, $false, $false, $false, $false, $true, 1, $false, "32-27=", 2) | Out-Null
$d.Content.Find.Execute("7+79=", $true, $false, $false, $false, $false, $true, 1, $false, "4+79=", 2) | Out-Null
$d.Content.Find.Execute("24+18=", $true, $false, $false, $false, $false, $true, 1, $false, "65+16=", 2) | Out-Null
$d.Content.Find.Execute("29+63=", $true, $false, $false, $false, $false, $true, 1, $false, "16+8=", 2) | Out-Null
$d.Content.Find.Execute("75-8=", $true, $false, $false, $false, $false, $true, 1, $false, "90-66=", 2) | Out-Null
$d.Content.Find.Execute("35+16=", $true, $false, $false, $false, $false, $true, 1, $false, "85-29=", 2) | Out-Null
$d.Content.Find.Execute("90-31=", $true, $false, $false, $false, $false, $true, 1, $false, "57+39=", 2) | Out-Null
$d.Content.Find.Execute("40-34=", $true, $false, $false, $false, $false, $true, 1, $false, "64-27=", 2) | Out-Null
$d.Content.Find.Execute("96-28=", $true, $false, $false, $false, $false, $true, 1, $false, "59+28=", 2) | Out-Null
$d.Content.Find.Execute("14+8=", $true, $false, $false, $false, $false, $true, 1, $false, "36-8=", 2) | Out-Null
$d.Content.Find.Execute("92-6=", $true, $false, $false, $false, $false, $true, 1, $false, "32-15=", 2) | Out-Null
$d.Content.Find.Execute("45+18=", $true, $false, $false, $false, $false, $true, 1, $false, "51-34=", 2) | Out-Null
$d.Content.Find.Execute("87-29=", $true, $false, $false, $false, $false, $true, 1, $false, "47+29=", 2) | Out-Null
$d.Content.Find.Execute("83-45=", $true, $false, $false, $false, $false, $true, 1, $false, "7+16=", 2) | Out-Null
$d.Content.Find.Execute("24+49=", $true, $false, $false, $false, $false, $true, 1, $false, "12+9=", 2) | Out-Null
$d.Content.Find.Execute("40-7=", $true, $false, $false, $false, $false, $true, 1, $false, "27+55=", 2) | Out-Null
$d.Content.Find.Execute("85-18=", $true, $false, $false, $false, $false, $true, 1, $false, "19+18=", 2) | Out-Null
$d.Content.Find.Execute("78+15=", $true, $false, $false, $false, $false, $true, 1, $false, "53-28=", 2) | Out-Null
$d.Content.Find.Execute("38+8=", $true, $false, $false, $false, $false, $true, 1, $false, "89+9=", 2) | Out-Null
$d.Content.Find.Execute("18+54=", $true, $false, $false, $false, $false, $true, 1, $false, "60-57=", 2) | Out-Null
$d.Content.Find.Execute("60-38=", $true, $false, $false, $false, $false, $true, 1, $false, "94-65=", 2) | Out-Null
$d.Content.Find.Execute("44-8=", $true, $false, $false, $false, $false, $true, 1, $false, "14+57=", 2) | Out-Null
$d.Content.Find.Execute("18+59=", $true, $false, $false, $false, $false, $true, 1, $false, "86-39=", 2) | Out-Null
$d.Content.Find.Execute("63-27=", $true, $false, $false, $false, $false, $true, 1, $false, "64-46=", 2) | Out-Null
$d.Content.Find.Execute("37+58=", $true, $false, $false, $false, $false, $true, 1, $false, "85-17=", 2) | Out-Null
$d.Content.Find.Execute("14-5=", $true, $false, $false, $false, $false, $true, 1, $false, "39+55=", 2) | Out-Null
$d.Content.Find.Execute("91-26=", $true, $false, $false, $false, $false, $true, 1, $false, "12-7=", 2) | Out-Null
$d.Content.Find.Execute("14+18=", $true, $false, $false, $false, $false, $true, 1, $false, "61-13=", 2) | Out-Null
$d.Content.Find.Execute("22-13=", $true, $false, $false, $false, $false, $true, 1, $false, "43-14=", 2) | Out-Null
$d.Content.Find.Execute("86-28=", $true, $false, $false, $false, $false, $true, 1, $false, "70-38=", 2) | Out-Null
$d.Content.Find.Execute("83-47=", $true, $false, $false, $false, $false, $true, 1, $false, "52-9=", 2) | Out-Null
$d.Content.Find.Execute("74-28=", $true, $false, $false, $false, $false, $true, 1, $false, "91-76=", 2) | Out-Null
$d.Content.Find.Execute("40-18=", $true, $false, $false, $false, $false, $true, 1, $false, "73-17=", 2) | Out-Null
$d.Content.Find.Execute("18+58=", $true, $false, $false, $false, $false, $true, 1, $false, "5+69=", 2) | Out-Null
$d.Content.Find.Execute("37+17=", $true, $false, $false, $false, $false, $true, 1, $false, "30-5=", 2) | Out-Null
$d.Content.Find.Execute("17+15=", $true, $false, $false, $false, $false, $true, 1, $false, "9+55=", 2) | Out-Null
$d.Content.Find.Execute("64+27=", $true, $false, $false, $false, $false, $true, 1, $false, "24-19=", 2) | Out-Null
$d.Content.Find.Execute("61-47=", $true, $false, $false, $false, $false, $true, 1, $false, "73-7=", 2) | Out-Null
$d.Content.Find.Execute("9+3=", $true, $false, $false, $false, $false, $true, 1, $false, "19+19=", 2) | Out-Null
$d.Content.Find.Execute("94-58=", $true, $false, $false, $false, $false, $true, 1, $false, "68-49=", 2) | Out-Null
$d.Content.Find.Execute("40-13=", $true, $false, $false, $false, $false, $true, 1, $false, "8+76=", 2) | Out-Null
$d.Content.Find.Execute("18+77=", $true, $false, $false, $false, $false, $true, 1, $false, "14+78=", 2) | Out-Null
$d.Content.Find.Execute("83+8=", $true, $false, $false, $false, $false, $true, 1, $false, "42+9=", 2) | Out-Null
$d.Content.Find.Execute("90-13=", $true, $false, $false, $false, $false, $true, 1, $false, "18+47=", 2) | Out-Null
$d.Content.Find.Execute("81-18=", $true, $false, $false, $false, $false, $true, 1, $false, "57+6=", 2) | Out-Null
$d.Content.Find.Execute("24+27=", $true, $false, $false, $false, $false, $true, 1, $false, "65-48=", 2) | Out-Null
$d.Content.Find.Execute("69+12=", $true, $false, $false, $false, $false, $true, 1, $false, "74-59=", 2) | Out-Null
$d.Content.Find.Execute("67+5=", $true, $false, $false, $false, $false, $true, 1, $false, "2+89=", 2) | Out-Null
$d.Content.Find.Execute("26+28=", $true, $false, $false, $false, $false, $true, 1, $false, "64+9=", 2) | Out-Null
$d.Content.Find.Execute("55+7=", $true, $false, $false, $false, $false, $true, 1, $false, "17+55=", 2) | Out-Null
$d.Content.Find.Execute("80-19=", $true, $false, $false, $false, $false, $true, 1, $false, "50-37=", 2) | Out-Null
$d.Content.Find.Execute("38+7=", $true, $false, $false, $false, $false, $true, 1, $false, "34-15=", 2) | Out-Null
$d.Content.Find.Execute("35+39=", $true, $false, $false, $false, $false, $true, 1, $false, "3+28=", 2) | Out-Null
$d.Content.Find.Execute("80-3=", $true, $false, $false, $false, $false, $true, 1, $false, "50-47=", 2) | Out-Null
$d.Content.Find.Execute("64-9=", $true, $false, $false, $false, $false, $true, 1, $false, "92-18=", 2) | Out-Null
$d.Content.Find.Execute("59+33=", $true, $false, $false, $false, $false, $true, 1, $false, "19+7=", 2) | Out-Null
$d.Content.Find.Execute("66-37=", $true, $false, $false, $false, $false, $true, 1, $false, "70-16=", 2) | Out-Null
